# Nexial "PlanTest-Script1.xlsx" style update:
# Add new JSON function `storeKeys(json,jsonpath,var)` to the "json" lookup
# list on the hidden '#system' sheet (inserted alphabetically at M16,
# pushing storeValue/storeValues down by one row), and remove the old
# standalone "text" lookup column (previously column Y), shifting the
# web/webalert/webcookie/ws/ws.async/xml lookup columns one letter to the
# left (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD) as well as removing the
# now-duplicated "text" category entry from the "target" category list in
# column A (row 25), shifting rows 26-31 up by one.

$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Column A ("target" category list): delete the "text" entry at A25,
#    shifting A26:A31 up to A25:A30 (A31 becomes blank / removed).
# ---------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $nextVal = $sys.Range("A" + ($r + 1)).Value()
    $sys.Range("A" + $r).Value = $nextVal
}
$sys.Range("A31").Value = ""

# ---------------------------------------------------------------------
# 2) Columns Y..AE: drop the old "text" column (Y), shifting every
#    following column left by one letter, for every row that holds data
#    (header row 1 plus data rows up to 129).
# ---------------------------------------------------------------------
$srcCols = @("Z", "AA", "AB", "AC", "AD", "AE")
$dstCols = @("Y", "Z", "AA", "AB", "AC", "AD")

for ($r = 1; $r -le 129; $r++) {
    for ($i = 0; $i -lt $srcCols.Length; $i++) {
        $srcVal = $sys.Range($srcCols[$i] + $r).Value()
        $sys.Range($dstCols[$i] + $r).Value = $srcVal
    }
    $sys.Range("AE" + $r).Value = ""
}

# ---------------------------------------------------------------------
# 3) Column M ("json" function list): insert the new `storeKeys` entry at
#    M16 (alphabetically between storeCount and storeValue), pushing the
#    existing storeValue/storeValues rows down to M17/M18.
# ---------------------------------------------------------------------
$sys.Range("M18").Value = $sys.Range("M17").Value()
$sys.Range("M17").Value = $sys.Range("M16").Value()
$sys.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 4) Update the affected defined names so they keep pointing at the
#    correct (now shifted) ranges.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
